$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." text. This is
# more robust than hard-coding a paragraph index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Ver no Jupiter*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $verPara = $d.Paragraphs.Item($targetIndex)

    # The paragraph immediately before it is the blank spacer paragraph that
    # was inserted along with the footer block, and the paragraph right
    # after it is the "(c) 2020 ... Creative Commons Attribution" footer
    # line. All three paragraphs (blank spacer, "Ver no Jupiter ...", and
    # the copyright line) were removed together, leaving the preceding
    # "Apostila ..." paragraph directly followed by the original trailing
    # blank paragraph / page break.
    $emptyPara = $verPara.Previous(1)
    $copyrightPara = $verPara.Next(1)

    $deleteRange = $d.Range($emptyPara.Range.Start, $copyrightPara.Range.End)
    $deleteRange.Delete()
}
